# Rename the "analysis" sheet to "data"
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("analysis")
$dataSheet.Name = "data"

# Re-activate the "data" sheet (it becomes the active tab again,
# moving the active tab away from "keywordsCleaned") and restore the
# previous selection to B53.
$dataSheet.Activate()
$dataSheet.Range("B53").Select()
